# "fixed dynamic browser implementation"
#
# Adds a per-sheet "browser" column to the ADDTOCART and BOOKFLIGHT test
# data sheets (ADDTOCART always runs on chrome, BOOKFLIGHT always runs on
# edge) mirroring the browser column that already existed on TESTRUNNER.
# Also updates the stored selections / active sheet to reflect where the
# author left off, and widens BOOKFLIGHT's first two columns to fit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ADDTOCART
$ws2 = $wb.Worksheets.Item(2)   # BOOKFLIGHT
$ws3 = $wb.Worksheets.Item(3)   # TESTRUNNER

# --- ADDTOCART: new column C = "browser", every row driven by chrome ---
$ws1.Range("C1").Value = "browser"
$ws1.Range("C2").Value = "chrome"
$ws1.Range("C3").Value = "chrome"
$ws1.Range("C4").Value = "chrome"
$ws1.Range("C5").Value = "chrome"

# --- BOOKFLIGHT: new column C = "browser", every row driven by edge ---
$ws2.Range("C1").Value = "browser"
$ws2.Range("C2").Value = "edge"
$ws2.Range("C3").Value = "edge"
$ws2.Range("C4").Value = "edge"

# Widen BOOKFLIGHT's departcity/arrivalcity columns now that the sheet
# carries more data.
$ws2.Columns.Item(1).ColumnWidth = 26 + 1/3
$ws2.Columns.Item(2).ColumnWidth = 23 + 2/3

# --- Restore the selections left on each sheet ---
$ws2.Range("D10").Select()
$ws3.Range("A5").Select()

# ADDTOCART is the sheet left active/visible when the file was saved.
$ws1.Activate()
$ws1.Range("E6").Select()
